# metodos para eliminar y para sumar las ventas y compras
# Append 32 additional rows (37-68) of data, following the same pattern
# as the existing rows (2-36): columns A-G filled with "jhon".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 37
$endRow = 68
$lastCol = 7  # columns A..G

for ($r = $startRow; $r -le $endRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = "jhon"
    }
}
